$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header column F1 = "OSMO_DEF", matching the formatting of the
# existing header cells (e.g. E1: bold, centered, bordered)
$ws.Range("F1").Value = "OSMO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
